# Update the "想去人数" (want-to-go count) figures in the 展览 and 全部类型
# sheets to reflect the latest generated data (gh-pages output @ 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 108
$ws1.Range("F4").Value  = 635
$ws1.Range("F5").Value  = 190
$ws1.Range("F7").Value  = 9615
$ws1.Range("F8").Value  = 866
$ws1.Range("F10").Value = 1216
$ws1.Range("F11").Value = 2238
$ws1.Range("F12").Value = 158
$ws1.Range("F14").Value = 5
$ws1.Range("F15").Value = 23
$ws1.Range("F17").Value = 464
$ws1.Range("F20").Value = 1345

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 108
$ws4.Range("F5").Value  = 635
$ws4.Range("F6").Value  = 190
$ws4.Range("F8").Value  = 9615
$ws4.Range("F9").Value  = 866
$ws4.Range("F11").Value = 1216
$ws4.Range("F12").Value = 2238
$ws4.Range("F13").Value = 158
$ws4.Range("F15").Value = 5
$ws4.Range("F16").Value = 23
$ws4.Range("F18").Value = 464
$ws4.Range("F21").Value = 1345

$wb.Save()
